$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15944
$ws1.Range("F8").Value = 718
$ws1.Range("F9").Value = 15498
$ws1.Range("F10").Value = 62
$ws1.Range("F11").Value = 9118
$ws1.Range("F12").Value = 399
$ws1.Range("F13").Value = 9
$ws1.Range("F14").Value = 1021
$ws1.Range("F20").Value = 69
$ws1.Range("F24").Value = 66
$ws1.Range("F25").Value = 1124
$ws1.Range("F33").Value = 68
$ws1.Range("F34").Value = 54
$ws1.Range("F37").Value = 462
$ws1.Range("F39").Value = 5611
$ws1.Range("F40").Value = 5234

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 76

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15944
$ws4.Range("F8").Value = 718
$ws4.Range("F9").Value = 15498
$ws4.Range("F10").Value = 62
$ws4.Range("F11").Value = 9118
$ws4.Range("F12").Value = 399
$ws4.Range("F13").Value = 9
$ws4.Range("F14").Value = 1021
$ws4.Range("F20").Value = 69
$ws4.Range("F24").Value = 66
$ws4.Range("F25").Value = 1124
$ws4.Range("F32").Value = 76
$ws4.Range("F35").Value = 68
$ws4.Range("F36").Value = 54
$ws4.Range("F39").Value = 462
$ws4.Range("F41").Value = 5611
$ws4.Range("F43").Value = 5234
